$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.436.24"
$ws.Range("E2").Value = "  +6.82%  "
$ws.Range("D3").Value = "3.537.12"
$ws.Range("E3").Value = "  +9.11%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  +10.25%  "
$ws.Range("D6").Value = "'559.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.81%  "
$ws.Range("D7").Value = "3.530.13"
$ws.Range("E7").Value = "  +8.78%  "
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'0.645"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.09%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.64%  "
$ws.Range("D12").Value = "'56.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.06%  "
$ws.Range("D13").Value = "'0.0000274"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.35%  "
$ws.Range("D14").Value = "'9.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.95%  "
$ws.Range("D15").Value = "4.110.12"
$ws.Range("E15").Value = "  +9.80%  "
$ws.Range("D16").Value = "3.544.21"
$ws.Range("E16").Value = "  +9.74%  "
$ws.Range("E17").Value = "  +5.66%  "
$ws.Range("D18").Value = "67.515.72"
$ws.Range("E18").Value = "  +7.23%  "
$ws.Range("D19").Value = "'18.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.40%  "
$ws.Range("D20").Value = "'11.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.66%  "
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "'408.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.31%  "
$ws.Range("E23").Value = "  +7.04%  "
$ws.Range("D24").Value = "'85.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.90%  "
$ws.Range("D25").Value = "'4.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.84%  "
$ws.Range("D26").Value = "'11.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'2.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.10%  "
$ws.Range("D28").Value = "'6.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  +7.56%  "
$ws.Range("D30").Value = "'8.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.49%  "
$ws.Range("D31").Value = "'30.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.22%  "
$ws.Range("D32").Value = "'684.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +7.94%  "
$ws.Range("D34").Value = "'11.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.71%  "
$ws.Range("E35").Value = "  +8.19%  "
$ws.Range("D36").Value = "'60.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("D37").Value = "'39.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.39%  "
$ws.Range("D38").Value = "0.0₃0826"
$ws.Range("E38").Value = "  +15.97%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'0.400"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.84%  "
$ws.Range("D41").Value = "'0.139"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.14%  "
$ws.Range("D42").Value = "'3.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.93%  "
$ws.Range("E43").Value = "  +18.21%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("E45").Value = "  +6.73%  "
$ws.Range("D46").Value = "3.030.96"
$ws.Range("E46").Value = "  +5.75%  "
$ws.Range("E47").Value = "  +8.06%  "
$ws.Range("D48").Value = "'3.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.15%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  +18.42%  "
$ws.Range("D51").Value = "'0.132"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.81%  "
